$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 29 (last quarter row) with refreshed metrics
$ws.Range("C29").Value = 213
$ws.Range("D29").Value = 34
$ws.Range("E29").Value = 179
$ws.Range("F29").Value = 5.851979345955249
